# Update the cryptos price table (columns D = Price, E = Volume(1h) change)
# with refreshed values, matching a GitHub-Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.024.94'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.233.42'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'293.57"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").Value = "'86.61"
$ws.Range("E6").Value = '  +5.34%  '
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").Value = "'30.78"
$ws.Range("E10").Value = '  +7.58%  '
$ws.Range("E11").Value = '  +2.10%  '
$ws.Range("D12").Value = "'47.08"
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("D14").Value = "'6.40"
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("D15").Value = '2.576.68'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '2.220.50'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = "'0.729"
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("D19").Value = '39.948.23'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  +2.77%  '
$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").Value = "'10.99"
$ws.Range("E22").Value = '  +9.14%  '
$ws.Range("D23").Value = "'65.27"
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").Value = "'235.59"
$ws.Range("E24").Value = '  +4.39%  '
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'2.47"
$ws.Range("E26").Value = '  +3.58%  '
$ws.Range("D27").Value = "'1.85"
$ws.Range("E27").Value = '  +5.68%  '
$ws.Range("D28").Value = "'22.81"
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("E29").Value = '  +2.88%  '
$ws.Range("D30").Value = "'9.27"
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").Value = "'33.13"
$ws.Range("E31").Value = '  +4.78%  '
$ws.Range("D32").Value = "'152.76"
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("D35").Value = "'0.0720"
$ws.Range("E35").Value = '  +3.63%  '
$ws.Range("E36").Value = '  +3.26%  '
$ws.Range("D37").Value = "'16.20"
$ws.Range("E37").Value = '  +9.88%  '
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").Value = "'0.0999"
$ws.Range("E40").Value = '  +5.62%  '
$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = '  +6.23%  '
$ws.Range("D42").Value = "'3.81"
$ws.Range("E42").Value = '  +4.30%  '
$ws.Range("D43").Value = '2.039.34'
$ws.Range("E43").Value = '  +7.23%  '
$ws.Range("E44").Value = '  +8.20%  '
$ws.Range("D45").Value = "'0.0271"
$ws.Range("E45").Value = '  +5.41%  '
$ws.Range("D46").Value = "'10.06"
$ws.Range("E46").Value = '  +11.91%  '
$ws.Range("D47").Value = "'16.65"
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("D48").Value = "'2.57"
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("D49").Value = '2.462.50'
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("D50").Value = "'70.94"
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("D51").Value = "'89.14"

Write-Host "Updated cryptos price table (rows 2-51)."
